$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume refresh for the cryptos table.
# Column D holds numeric-looking text (e.g. "1.001", "24.270.33") that must
# stay plain text, so we force NumberFormat "@" before writing those values.
# Columns B/C/E are unambiguous text and need no special handling.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.270.33"
$ws.Range("E2").Value = "  +15.07%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.678.76"
$ws.Range("E3").Value = "  +9.09%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.99%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.24"
$ws.Range("E5").Value = "  +9.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9965"
$ws.Range("E6").Value = "  +3.61%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3723"
$ws.Range("E7").Value = "  +2.88%  "

# Row 8
$ws.Range("E8").Value = "  +8.24%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.18"
$ws.Range("E9").Value = "  +18.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.187"
$ws.Range("E10").Value = "  +8.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07289"
$ws.Range("E11").Value = "  +7.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9986"
$ws.Range("E12").Value = "  -0.82%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.60"
$ws.Range("E13").Value = "  +10.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.104"
$ws.Range("E14").Value = "  +7.77%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.771"
$ws.Range("E15").Value = "  +6.69%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.677.90"
$ws.Range("E16").Value = "  +8.92%  "

# Row 17
$ws.Range("E17").Value = "  +6.65%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9968"
$ws.Range("E18").Value = "  +3.75%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06706"
$ws.Range("E19").Value = "  +10.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.74"
$ws.Range("E20").Value = "  +13.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.49"
$ws.Range("E21").Value = "  +10.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.126"
$ws.Range("E22").Value = "  +7.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.01"
$ws.Range("E23").Value = "  +5.97%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.232.91"
$ws.Range("E24").Value = "  +13.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.400"
$ws.Range("E25").Value = "  +3.45%  "

# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.358"
$ws.Range("E26").Value = "  -8.90%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.672"
$ws.Range("E27").Value = "  +21.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.34"
$ws.Range("E28").Value = "  +2.10%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.51"
$ws.Range("E29").Value = "  +10.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.861.79"
$ws.Range("E30").Value = "  +8.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.14"
$ws.Range("E31").Value = "  +7.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.422"
$ws.Range("E32").Value = "  +23.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.040"
$ws.Range("E33").Value = "  +0.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9941"
$ws.Range("E34").Value = "  +17.08%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.747"
$ws.Range("E35").Value = "  +15.99%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08445"
$ws.Range("E36").Value = "  +5.57%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.44"
$ws.Range("E37").Value = "  +16.68%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.370"
$ws.Range("E38").Value = "  +8.52%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06423"
$ws.Range("E39").Value = "  +9.66%  "

# Row 40
$ws.Range("E40").Value = "  +15.53%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.296"
$ws.Range("E41").Value = "  +7.32%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02347"
$ws.Range("E42").Value = "  +11.58%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2114"
$ws.Range("E43").Value = "  +10.56%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6160"
$ws.Range("E44").Value = "  +13.16%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9956"
$ws.Range("E45").Value = "  +3.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.801"
$ws.Range("E46").Value = "  +6.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.21"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5978"
$ws.Range("E48").Value = "  +9.83%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.55"
$ws.Range("E49").Value = "  +4.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.022"
$ws.Range("E50").Value = "  +8.35%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07141"
$ws.Range("E51").Value = "  +8.32%  "
